$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: 001 -> 002 (must remain text with leading zero, preserve the
# cell's original "no explicit style" state)
$origStyleJ2 = $ws.Range("J2").Style
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = $origStyleJ2

# REPORT_DATE: 2018-12-31 -> 2020-06-30 (stored as plain text, not a real date)
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures for the new (2020-06-30) reporting period
$ws.Range("O2").Value = 49097516.83
$ws.Range("P2").Value = 646129750.46
$ws.Range("Q2").Value = 558798267.46
$ws.Range("R2").Value = -27.8537473766
$ws.Range("S2").Value = 438334064.1
$ws.Range("T2").Value = 438334064.1
$ws.Range("U2").Value = -26.8166508065
$ws.Range("V2").Value = 60705688.34
$ws.Range("W2").Value = 25697899.9
$ws.Range("X2").Value = 8076502.46
$ws.Range("Y2").Value = 53565643.74
$ws.Range("Z2").Value = 55792894.12
$ws.Range("AA2").Value = 6860857.31
$ws.Range("AG2").Value = 7618470.97
$ws.Range("AP2").Value = -26.9309106621
$ws.Range("AQ2").Value = -36.436252726036
$ws.Range("AR2").Value = -27.879787061644
$ws.Range("AS2").Value = 42907748.31
$ws.Range("AT2").Value = -34.719031574241
